# Amend corrected label annotations
# Normalize the category tags in column F to lowercase (row 41 also has its
# two "||"-joined segments reordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16  = "93_referral_statement"
    17  = "93_referral_statement"
    29  = "18_hazards_to_humans_and_domestic_animals"
    31  = "ppe"
    32  = "ppe"
    37  = "application instructions"
    38  = "env warning - species"
    39  = "env warning - species"
    40  = "env warning - water"
    41  = "env warning - water || off target movement"
    51  = "pollinator"
    53  = "application instructions"
    54  = "application instructions"
    55  = "application instructions"
    56  = "application instructions"
    58  = "135_product_information"
    60  = "application instructions"
    62  = "use restrictions"
    67  = "irrigation"
    68  = "application instructions"
    69  = "chemigation"
    70  = "safety procedures"
    71  = "irrigation"
    72  = "irrigation"
    73  = "chemigation"
    74  = "irrigation"
    75  = "irrigation"
    76  = "safety procedures"
    77  = "safety procedures"
    78  = "off target movement"
    79  = "off target movement"
    80  = "off target movement"
    82  = "off target movement"
    84  = "off target movement"
    85  = "off target movement"
    87  = "off target movement"
    89  = "172_sensitive_areas"
    90  = "172_sensitive_areas"
    91  = "application instructions"
    92  = "mixing"
    93  = "mixing"
    95  = "mixing"
    96  = "mixing"
    97  = "use restrictions"
    282 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
